$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: B16 becomes a real number (was text "3"); C16 becomes an empty
# text cell (was text "nan"). Using a lone apostrophe forces a text
# (quote-prefixed) entry whose displayed/stored text is empty, then the
# style is reset to Normal so no quote-prefix formatting lingers.
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = "'"
$ws.Range("C16").Style = "Normal"

# New row 17 data (B17 must stay text "4", not become a number, so it is
# entered with a quote-prefix and then the quote-prefix style is cleared)
$ws.Range("A17").Value = "parisk"
$ws.Range("B17").Value = "'4"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "show through several experiments"
$ws.Range("D17").Value = "FBK"
$ws.Range("E17").Value = "RES"
$ws.Range("F17").Value = "9cb2103f-10a8-4188-b35f-b6e342d90889"
$ws.Range("G17").Value = "rJwelMbR-_annotated.xlsx"
$ws.Range("H17").Value = "The authors show through several experiments that the divide and conquer (DnC) technique can solve more complex tasks than can be solved with conventional policy gradient methods (TRPO is used as the baseline)."
$ws.Range("I17").Value = "Correct"
